# Saldo.xlsx update:
#  - Remove 6 stale/duplicate balance rows (by account number).
#  - Move the Priscilla row (004224284) further down the list (after Silvio /
#    008197302) and correct her balance from 82606.04 to 606.04.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-to-top so earlier row numbers stay valid as we go.
# (Row numbers below are the ORIGINAL 1-based row numbers in the sheet.)
#   Row 20 -> 004853111 Marcondes 1000
#   Row 19 -> 004498637 Tiago     1000
#   Row 15 -> 005624730 Isabel    7000
#   Row 12 -> 005002457 Rosangela 14043.01
#   Row 10 -> 004515341 Bruno     25083.33
#   Row 9  -> 004224815 Guilherme 35404.38
#   Row 4  -> 004224284 Priscilla 82606.04 (will be re-inserted below, fixed)
$rowsToDelete = @(20, 19, 15, 12, 10, 9, 4)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# After those 7 deletions, "008197302 Silvio" (originally row 47) now sits at
# row 40. Insert a fresh row right after it (row 41) for Priscilla, with the
# corrected balance.
$ws.Rows.Item(41).Insert()
$ws.Cells.Item(41, 1).Value = "'004224284"
$ws.Cells.Item(41, 2).Value = "Priscilla"
$ws.Cells.Item(41, 3).Value = 606.04
